$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 175.16667
$ws.Range("I6").Value = 134.2
$ws.Range("J6").Value = 380
$ws.Range("K6").Value = 402.6
$ws.Range("L6").Value = 1140
$ws.Range("M6").Value = -290.6
$ws.Range("N6").Value = -1364
$ws.Range("H80").Value = 3468.2778
$ws.Range("I80").Value = 802
$ws.Range("J80").Value = 6134.5557
$ws.Range("K80").Value = 2406
$ws.Range("L80").Value = 18403.6671
$ws.Range("M80").Value = -1408
$ws.Range("N80").Value = -20399.6671
$ws.Range("H83").Value = 3468.2778
$ws.Range("I83").Value = 802
$ws.Range("J83").Value = 6134.5557
$ws.Range("K83").Value = 7218
$ws.Range("L83").Value = 55211.0013
$ws.Range("M83").Value = -2226
$ws.Range("N83").Value = -65195.0013
$ws.Range("H98").Value = 1115.0646
$ws.Range("I98").Value = 1122.9434
$ws.Range("K98").Value = 1122.9434
$ws.Range("M98").Value = 375.0565999999999
$ws.Range("H122").Value = 1115.0646
$ws.Range("I122").Value = 1122.9434
$ws.Range("K122").Value = 3368.8302
$ws.Range("M122").Value = -918.8302000000003
$ws.Range("H137").Value = 39749.914
$ws.Range("I137").Value = 49499.324
$ws.Range("J137").Value = 3677.1
$ws.Range("K137").Value = 148497.972
$ws.Range("L137").Value = 11031.3
$ws.Range("M137").Value = -145947.972
$ws.Range("N137").Value = -16131.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 513.1818
$ws.Range("I4").Value = 406.2857
$ws.Range("J4").Value = 700.25
$ws.Range("K4").Value = 406.2857
$ws.Range("L4").Value = 700.25
$ws.Range("M4").Value = -290.2857
$ws.Range("N4").Value = -932.25
$ws.Range("H14").Value = 560.6
$ws.Range("I14").Value = 435
$ws.Range("J14").Value = 749
$ws.Range("K14").Value = 435
$ws.Range("L14").Value = 749
$ws.Range("M14").Value = -260
$ws.Range("N14").Value = -1099
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H32").Value = 8647.886
$ws.Range("I32").Value = 4609.8135
$ws.Range("K32").Value = 4609.8135
$ws.Range("M32").Value = -4322.8135
$ws.Range("H33").Value = 18642.857
$ws.Range("J33").Value = 19333.334
$ws.Range("L33").Value = 19333.334
$ws.Range("N33").Value = -19991.334
$ws.Range("H74").Value = 24359.55
$ws.Range("I74").Value = 1271.2858
$ws.Range("K74").Value = 1271.2858
$ws.Range("M74").Value = -397.2858000000001
$ws.Range("H77").Value = 24359.55
$ws.Range("I77").Value = 1271.2858
$ws.Range("K77").Value = 6356.429
$ws.Range("M77").Value = -1988.429
$ws.Range("H88").Value = 1215.8235
$ws.Range("I88").Value = 1670.8334
$ws.Range("J88").Value = 967.63635
$ws.Range("K88").Value = 1670.8334
$ws.Range("L88").Value = 967.63635
$ws.Range("M88").Value = -1264.8334
$ws.Range("N88").Value = -1779.63635
$ws.Range("H91").Value = 1215.8235
$ws.Range("I91").Value = 1670.8334
$ws.Range("J91").Value = 967.63635
$ws.Range("K91").Value = 1670.8334
$ws.Range("L91").Value = 967.63635
$ws.Range("M91").Value = -266.8334
$ws.Range("N91").Value = -3775.63635
$ws.Range("H132").Value = 2232.6086
$ws.Range("I132").Value = 1806.6511
$ws.Range("K132").Value = 5419.9533
$ws.Range("M132").Value = -2889.9533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 297.66666
$ws.Range("I12").Value = 196.5
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 196.5
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -28.5
$ws.Range("N12").Value = -836
$ws.Range("H22").Value = 3949.6667
$ws.Range("I22").Value = 6066
$ws.Range("K22").Value = 6066
$ws.Range("M22").Value = -5893
$ws.Range("H134").Value = 3724.2922
$ws.Range("I134").Value = 2055.721
$ws.Range("K134").Value = 6167.163
$ws.Range("M134").Value = -3632.163

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 635.7778
$ws.Range("I7").Value = 431.75
$ws.Range("K7").Value = 431.75
$ws.Range("M7").Value = -318.75
$ws.Range("H16").Value = 1836.5
$ws.Range("I16").Value = 1341
$ws.Range("K16").Value = 1341
$ws.Range("M16").Value = -1054
$ws.Range("H31").Value = 24160.617
$ws.Range("I31").Value = 3786.1667
$ws.Range("J31").Value = 27142.244
$ws.Range("K31").Value = 3786.1667
$ws.Range("L31").Value = 27142.244
$ws.Range("M31").Value = -3491.1667
$ws.Range("N31").Value = -27732.244
$ws.Range("H34").Value = 24160.617
$ws.Range("I34").Value = 3786.1667
$ws.Range("J34").Value = 27142.244
$ws.Range("K34").Value = 3786.1667
$ws.Range("L34").Value = 27142.244
$ws.Range("M34").Value = -3584.1667
$ws.Range("N34").Value = -27546.244
$ws.Range("H58").Value = 4666.5386
$ws.Range("I58").Value = 5467.4
$ws.Range("J58").Value = 3574.4546
$ws.Range("K58").Value = 5467.4
$ws.Range("L58").Value = 3574.4546
$ws.Range("M58").Value = -5264.4
$ws.Range("N58").Value = -3980.4546
$ws.Range("H107").Value = 2588.5881
$ws.Range("I107").Value = 2071.1428
$ws.Range("J107").Value = 5003.3335
$ws.Range("K107").Value = 2071.1428
$ws.Range("L107").Value = 5003.3335
$ws.Range("M107").Value = -151.1428000000001
$ws.Range("N107").Value = -8843.333500000001
$ws.Range("H113").Value = 1836.5
$ws.Range("I113").Value = 1341
$ws.Range("K113").Value = 1341
$ws.Range("M113").Value = 829
$ws.Range("H121").Value = 28279.8
$ws.Range("J121").Value = 28279.8
$ws.Range("L121").Value = 28279.8
$ws.Range("N121").Value = -30899.8
$ws.Range("H132").Value = 45912.316
$ws.Range("I132").Value = 29493.459
$ws.Range("J132").Value = 132697.72
$ws.Range("K132").Value = 88480.37699999999
$ws.Range("L132").Value = 398093.16
$ws.Range("M132").Value = -85950.37699999999
$ws.Range("N132").Value = -403153.16
$ws.Range("H136").Value = 4666.5386
$ws.Range("I136").Value = 5467.4
$ws.Range("J136").Value = 3574.4546
$ws.Range("K136").Value = 16402.2
$ws.Range("L136").Value = 10723.3638
$ws.Range("M136").Value = -13852.2
$ws.Range("N136").Value = -15823.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 99.15385000000001
$ws.Range("J34").Value = 100
$ws.Range("L34").Value = 300
$ws.Range("N34").Value = -468
$ws.Range("H55").Value = 93198.91
$ws.Range("I55").Value = 797.25
$ws.Range("J55").Value = 145999.86
$ws.Range("K55").Value = 2391.75
$ws.Range("L55").Value = 437999.58
$ws.Range("M55").Value = -2214.75
$ws.Range("N55").Value = -438353.58
$ws.Range("H131").Value = 9061254
$ws.Range("J131").Value = 9262747
$ws.Range("L131").Value = 27788241
$ws.Range("N131").Value = -27798321

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 25644084
$ws.Range("J2").Value = 58829800
$ws.Range("L2").Value = 58829800
$ws.Range("N2").Value = -58830026

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6778.385
$ws.Range("I7").Value = 2783.1667
$ws.Range("J7").Value = 10202.857
$ws.Range("K7").Value = 2783.1667
$ws.Range("L7").Value = 10202.857
$ws.Range("M7").Value = -2671.1667
$ws.Range("N7").Value = -10426.857
$ws.Range("H29").Value = 24999
$ws.Range("I29").Value = 24999
$ws.Range("K29").Value = 24999
$ws.Range("M29").Value = -24704
$ws.Range("H40").Value = 13334
$ws.Range("I40").Value = 10004
$ws.Range("K40").Value = 10004
$ws.Range("M40").Value = -9868
$ws.Range("H42").Value = 9675
$ws.Range("I42").Value = 9675
$ws.Range("K42").Value = 9675
$ws.Range("M42").Value = -9112
$ws.Range("H49").Value = 9675
$ws.Range("I49").Value = 9675
$ws.Range("K49").Value = 9675
$ws.Range("M49").Value = -9528
$ws.Range("H126").Value = 6778.385
$ws.Range("I126").Value = 2783.1667
$ws.Range("J126").Value = 10202.857
$ws.Range("K126").Value = 8349.500100000001
$ws.Range("L126").Value = 30608.571
$ws.Range("M126").Value = -5879.500100000001
$ws.Range("N126").Value = -35548.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""
$ws.Range("H38").Value = 156
$ws.Range("I38").Value = 156
$ws.Range("K38").Value = 156
$ws.Range("M38").Value = 317
$ws.Range("H51").Value = 20373.6
$ws.Range("I51").Value = 12967.25
$ws.Range("K51").Value = 12967.25
$ws.Range("M51").Value = -12457.25
$ws.Range("H96").Value = 3228.5652
$ws.Range("I96").Value = 3110.0625
$ws.Range("K96").Value = 3110.0625
$ws.Range("M96").Value = -1737.0625
$ws.Range("H122").Value = 3212.6667
$ws.Range("I122").Value = 1517.4445
$ws.Range("K122").Value = 4552.333500000001
$ws.Range("M122").Value = -2102.333500000001
